$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy the header style (bold, centered, bordered) from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..38: I = 1 (except row 36 = 4), J = H value (except row 36 = 6)
for ($r = 2; $r -le 38; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    if ($r -eq 36) {
        $ws.Cells.Item($r, 9).Value = 4
        $ws.Cells.Item($r, 10).Value = 6
    } else {
        $ws.Cells.Item($r, 9).Value = 1
        $ws.Cells.Item($r, 10).Value = $hVal
    }
}
